# Rename sheets (old LogParser/LTE naming -> new naming, per commit 'modify logic for reading logs')
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("LTE_Power_Current").Name = "Zigbee_Power_Current"
$wb.Worksheets.Item("LTE_LAN").Name = "Zigbee_LAN"
$wb.Worksheets.Item("Zigbee_Current").Name = "LTE_Current"
$wb.Worksheets.Item("Zigbee_dBm").Name = "LTE_dBm"

# Rebuild each sheet's data rows (rows 2-10) sorted by SN (serial number),
# and with corrected Current/dBm values for the re-parsed log that was previously
# mis-attributed to SN 5891801300003_0729112049.

# --- Zigbee_Power_Current ---
$ws = $wb.Worksheets.Item("Zigbee_Power_Current")
$data = @(
    @("5891801300001_0731161805", 9.442, 9.688, 8.916, 81, 82, 80),
    @("5891801300002_0731153739", 8.998, 8.928, 9.047, 81, 79, 79),
    @("5891801300003_0729112049", 8.722999999999999, 8.989, 8.476, 79, 77, 78),
    @("5891801300004_0731153126", 9.665999999999999, 8.817, 9.260000000000002, 81, 81, 80),
    @("5891801300005_0729063540", 9.238, 8.575000000000001, 8.973, 77, 77, 76),
    @("5891801300006_0731155333", 8.623, 8.924, 9.232000000000001, 79, 73, 81),
    @("5891801300008_0727181606", "", "", "", "", "", ""),
    @("5891801300008_0728214258", 9.258, 9.426, 8.617, 81, 83, 82),
    @("5891801300009_0729031917", 8.575999999999999, 8.688, 8.975000000000001, 78, 77, 76)
)
$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ("$val" -eq "") {
            $ws.Cells.Item($r, $c + 1).ClearContents()
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
    $r++
}

# --- Zigbee_LAN ---
$ws = $wb.Worksheets.Item("Zigbee_LAN")
$data = @(
    @("5891801300001_0731161805", -15, -28),
    @("5891801300002_0731153739", -16, -28),
    @("5891801300003_0729112049", -15, -28),
    @("5891801300004_0731153126", -16, -29),
    @("5891801300005_0729063540", -15, -28),
    @("5891801300006_0731155333", -16, -28),
    @("5891801300008_0727181606", "", ""),
    @("5891801300008_0728214258", -13, -26),
    @("5891801300009_0729031917", -16, -29)
)
$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ("$val" -eq "") {
            $ws.Cells.Item($r, $c + 1).ClearContents()
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
    $r++
}

# --- LTE_Current ---
$ws = $wb.Worksheets.Item("LTE_Current")
$data = @(
    @("5891801300001_0731161805", 246, 213, 156),
    @("5891801300002_0731153739", 252, 211, 153),
    @("5891801300003_0729112049", 244, 204, 152),
    @("5891801300004_0731153126", 241, 205, 156),
    @("5891801300005_0729063540", 251, 213, 155),
    @("5891801300006_0731155333", 241, 207, 155),
    @("5891801300008_0727181606", 245, 209, 153),
    @("5891801300008_0728214258", 243, 209, 154),
    @("5891801300009_0729031917", 243, 213, 135)
)
$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ("$val" -eq "") {
            $ws.Cells.Item($r, $c + 1).ClearContents()
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
    $r++
}

# --- LTE_dBm ---
$ws = $wb.Worksheets.Item("LTE_dBm")
$data = @(
    @("5891801300001_0731161805", 22.592, 21.27, 31.718, -58),
    @("5891801300002_0731153739", 23.126, 21.759, 31.393, -58),
    @("5891801300003_0729112049", 23.071, 20.848, 31.622, -588),
    @("5891801300004_0731153126", 22.453, 21.119, 31.638, -58),
    @("5891801300005_0729063540", 23.045, 20.841, 31.558, -58),
    @("5891801300006_0731155333", 22.564, 21.071, 31.667, -58),
    @("5891801300008_0727181606", 21.654, 20.876, 30.265, -58),
    @("5891801300008_0728214258", 23.047, 20.811, 31.579, -58),
    @("5891801300009_0729031917", 22.954, 21.161, 31.305, -58)
)
$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ("$val" -eq "") {
            $ws.Cells.Item($r, $c + 1).ClearContents()
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
    $r++
}
